$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 6046
$ws.Range('J3').Value = 6467
$ws.Range('J4').Value = 1397
$ws.Range('J5').Value = 494
$ws.Range('J6').Value = 8371
$ws.Range('J7').Value = 22775

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J2').Value = 180
$ws.Range('J4').Value = 99
$ws.Range('J5').Value = 72
$ws.Range('J6').Value = 171
$ws.Range('J7').Value = 672
$ws.Range('J8').Value = 1429
$ws.Range('J11').Value = 364
$ws.Range('J18').Value = 192
$ws.Range('J19').Value = 674
$ws.Range('J20').Value = 479
$ws.Range('J24').Value = 73
$ws.Range('J25').Value = 112
$ws.Range('J27').Value = 141
$ws.Range('J29').Value = 1256
$ws.Range('J30').Value = 86
$ws.Range('J31').Value = 206
$ws.Range('J33').Value = 1049
$ws.Range('J34').Value = 104
$ws.Range('J37').Value = 697
$ws.Range('J41').Value = 153
$ws.Range('J42').Value = 959
$ws.Range('J43').Value = 191
$ws.Range('J44').Value = 175
$ws.Range('J46').Value = 75
$ws.Range('J50').Value = 138
$ws.Range('J52').Value = 569
$ws.Range('J53').Value = 319
$ws.Range('J54').Value = 439
$ws.Range('J55').Value = 323
$ws.Range('J56').Value = 29
$ws.Range('J63').Value = 78
$ws.Range('J65').Value = 562
$ws.Range('J67').Value = 864
$ws.Range('J73').Value = 217
$ws.Range('J75').Value = 70
$ws.Range('J76').Value = 348
$ws.Range('J77').Value = 171
$ws.Range('J79').Value = 651
$ws.Range('J80').Value = 35
$ws.Range('J81').Value = 24
$ws.Range('J83').Value = 452
$ws.Range('J85').Value = 936
$ws.Range('J88').Value = 238
$ws.Range('J89').Value = 300
$ws.Range('J90').Value = 244
$ws.Range('J91').Value = 262
$ws.Range('J94').Value = 234
$ws.Range('J96').Value = 258
$ws.Range('J98').Value = 167
$ws.Range('J99').Value = 354
$ws.Range('J101').Value = 22775

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J3').Value = 71
$ws.Range('J7').Value = 258

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J6').Value = 217
$ws.Range('J7').Value = 672

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J3').Value = 70
$ws.Range('J6').Value = 156
$ws.Range('J7').Value = 364

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J2').Value = 93
$ws.Range('J3').Value = 85
$ws.Range('J7').Value = 300

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 249
$ws.Range('J7').Value = 936

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J3').Value = 173
$ws.Range('J4').Value = 21
$ws.Range('J5').Value = 10
$ws.Range('J7').Value = 569

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J4').Value = 11
$ws.Range('J6').Value = 210
$ws.Range('J7').Value = 319

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 388
$ws.Range('J3').Value = 435
$ws.Range('J4').Value = 78
$ws.Range('J7').Value = 1429

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J3').Value = 168
$ws.Range('J7').Value = 452

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 247
$ws.Range('J3').Value = 349
$ws.Range('J6').Value = 366
$ws.Range('J7').Value = 1049

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J3').Value = 240
$ws.Range('J7').Value = 697

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 164
$ws.Range('J7').Value = 562

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J3').Value = 139
$ws.Range('J7').Value = 354

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('J2').Value = 30
$ws.Range('J7').Value = 86

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J3').Value = 56
$ws.Range('J7').Value = 206

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J3').Value = 327
$ws.Range('J7').Value = 864

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J6').Value = 210
$ws.Range('J7').Value = 439

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 385
$ws.Range('J3').Value = 434
$ws.Range('J4').Value = 68
$ws.Range('J6').Value = 321
$ws.Range('J7').Value = 1256

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 162
$ws.Range('J3').Value = 197
$ws.Range('J7').Value = 674

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J2').Value = 55
$ws.Range('J6').Value = 71
$ws.Range('J7').Value = 175

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J2').Value = 57
$ws.Range('J4').Value = 27
$ws.Range('J6').Value = 193
$ws.Range('J7').Value = 348

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J2').Value = 51
$ws.Range('J6').Value = 63
$ws.Range('J7').Value = 171

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J2').Value = 33
$ws.Range('J3').Value = 22
$ws.Range('J7').Value = 153

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J6').Value = 500
$ws.Range('J7').Value = 959

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J3').Value = 69
$ws.Range('J6').Value = 169
$ws.Range('J7').Value = 323

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('J4').Value = 12
$ws.Range('J6').Value = 18
$ws.Range('J7').Value = 73

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('J6').Value = 30
$ws.Range('J7').Value = 75

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J2').Value = 73
$ws.Range('J6').Value = 63
$ws.Range('J7').Value = 262

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J2').Value = 181
$ws.Range('J3').Value = 223
$ws.Range('J4').Value = 39
$ws.Range('J6').Value = 190
$ws.Range('J7').Value = 651

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J3').Value = 166
$ws.Range('J4').Value = 42
$ws.Range('J7').Value = 479

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('J2').Value = 51
$ws.Range('J3').Value = 41
$ws.Range('J6').Value = 91
$ws.Range('J7').Value = 192

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('J6').Value = 38
$ws.Range('J7').Value = 104

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J2').Value = 43
$ws.Range('J6').Value = 130
$ws.Range('J7').Value = 234

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('J2').Value = 50
$ws.Range('J7').Value = 112

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('J6').Value = 103
$ws.Range('J7').Value = 167

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J3').Value = 35
$ws.Range('J4').Value = 22
$ws.Range('J6').Value = 44
$ws.Range('J7').Value = 138

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J6').Value = 70
$ws.Range('J7').Value = 217

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('J6').Value = 69
$ws.Range('J7').Value = 180

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J3').Value = 64
$ws.Range('J7').Value = 238

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('J6').Value = 34
$ws.Range('J7').Value = 72

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J3').Value = 34
$ws.Range('J7').Value = 141

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('J6').Value = 16
$ws.Range('J7').Value = 70

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J2').Value = 86
$ws.Range('J7').Value = 244

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('J6').Value = 114
$ws.Range('J7').Value = 191

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J6').Value = 29
$ws.Range('J7').Value = 171

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range('J6').Value = 16
$ws.Range('J7').Value = 29

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('J6').Value = 18
$ws.Range('J7').Value = 35

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('J2').Value = 33
$ws.Range('J7').Value = 99

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range('J6').Value = 7
$ws.Range('J7').Value = 24
